$p = $ppt.ActivePresentation

# --- Slide 7: "TextBox 4" ($ anaconda search  --channel r  gsl) ---
# Merge the "anaconda " run and the "search  --channel r  gsl" run into a
# single run (leaving the separate leading "$ " run untouched).
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(4)
$tr7 = $shp7.TextFrame.TextRange
$full7 = $tr7.Text
$idx7 = $full7.IndexOf("anaconda") + 1
$len7 = $full7.Length - $idx7 + 1
$merge7 = $tr7.Characters($idx7, $len7)
$merge7.Text = "anaconda search  --channel r  gsl"

# --- Slide 9: "TextBox 5" ($ conda env  export  R_base_environment.yml) ---
# Fix the typo: replace the double space before the filename with a
# redirection operator "> " so the command reads
# "$ conda env  export > R_base_environment.yml".
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(5)
$tr9 = $shp9.TextFrame.TextRange
$full9 = $tr9.Text
$marker9 = "export "
$idx9 = $full9.IndexOf($marker9) + $marker9.Length + 1
$fix9 = $tr9.Characters($idx9, 1)
$fix9.Text = "> "
